$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (defined name / autofilter formula follow automatically)
$ws.Name = "EPEX Spot Results"

# Shift the header row one column to the left (B:H -> A:G) and add a new
# "Volume" header in the newly vacated H1 cell.
$ws.Range("A1").Value = "HH"
$ws.Range("B1").Value = "Low"
$ws.Range("C1").Value = "High"
$ws.Range("D1").Value = "Last"
$ws.Range("E1").Value = "Weight Avg."
$ws.Range("F1").Value = "Buy Volume"
$ws.Range("G1").Value = "Sell Volume"
$ws.Range("H1").Value = "Volume"

# Re-apply the autofilter so it covers A1:H1 instead of B1:H1
$ws.AutoFilterMode = $false
$ws.Range("A1:H1").AutoFilter()

# The _FilterDatabase defined name doesn't follow the autofilter range
# automatically - point it at the new range (and new sheet name) explicitly.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "='EPEX Spot Results'!`$A`$1:`$H`$1"

# Column widths: A:D share one width, E:H share a (wider) width
$ws.Range("A1:D1").ColumnWidth = 15.69921875
$ws.Range("E1:H1").ColumnWidth = 22.69921875

# Move the active selection to A2 (was B2)
$ws.Range("A2").Select()

# Add the "Internal Only" footer
$ws.PageSetup.RightFooter = "`r&1#&`"Aptos`"&10&K000000 Internal Only"
